# Updated cryptos list (price/volume refresh) on Mon Oct 16 22:41:00 UTC 2023 with GitHub Actions
# For cells whose new value is a plain numeric-looking string, we briefly force the
# cell to Text format so Excel stores it as a string (matching the source data,
# which keeps these as text, e.g. "215.01") instead of auto-converting to a number
# and losing precision/formatting. ClearFormats() afterwards restores the cell to
# its original (default) style so no visible formatting changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.459.54'
$ws.Range('E2').Value = '  +4.22%  '
$ws.Range('D3').Value = '1.598.96'
$ws.Range('E3').Value = '  +2.03%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.01'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.498'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.53%  '
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.08'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +9.09%  '
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('D12').Value = '1.824.59'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '1.580.83'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('E14').Value = '  +0.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.535'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.15%  '
$ws.Range('D16').Value = '28.471.97'
$ws.Range('E16').Value = '  +4.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.35'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '232.89'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +7.07%  '
$ws.Range('D19').Value = '0.0₃0713'
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.54'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.42'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.09%  '
$ws.Range('E24').Value = '  +1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.88'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('E30').Value = '  +0.84%  '
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('E33').Value = '  +0.78%  '
$ws.Range('D34').Value = '1.422.68'
$ws.Range('E34').Value = '  -0.95%  '
$ws.Range('E35').Value = '  -0.54%  '
$ws.Range('E36').Value = '  -4.09%  '
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0168'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.53'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +8.20%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.545'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.823'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('E42').Value = '  -2.79%  '
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.84'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +6.20%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.980'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.94'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('D47').Value = '1.736.34'
$ws.Range('E47').Value = '  +1.96%  '
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('E50').Value = '  +5.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0526'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.12%  '
